$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.439.68"
$ws.Range("E2").Value = "  -1.90%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.490.28"
$ws.Range("E3").Value = "  -2.20%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "611.98"
$ws.Range("E5").Value = "  +5.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "189.20"
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.625"
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.212"
$ws.Range("E9").Value = "  -3.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.645"
$ws.Range("E10").Value = "  -1.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.79"
$ws.Range("E11").Value = "  -3.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000306"
$ws.Range("E12").Value = "  -3.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.46"
$ws.Range("E13").Value = "  -0.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.053.86"
$ws.Range("E14").Value = "  -1.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "601.53"
$ws.Range("E15").Value = "  +3.98%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.533.58"
$ws.Range("E16").Value = "  -1.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.94"
$ws.Range("E17").Value = "  -1.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.55"
$ws.Range("E18").Value = "  -2.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.511.08"
$ws.Range("E19").Value = "  -1.50%  "
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.983"
$ws.Range("E21").Value = "  -2.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.08"
$ws.Range("E22").Value = "  -3.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "105.91"
$ws.Range("E23").Value = "  +12.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.70"
$ws.Range("E24").Value = "  +2.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.13"
$ws.Range("E25").Value = "  +5.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.02"
$ws.Range("E26").Value = "  +1.64%  "
$ws.Range("E27").Value = "  -2.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.66"
$ws.Range("E28").Value = "  +3.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.31"
$ws.Range("E29").Value = "  +2.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.91"
$ws.Range("E30").Value = "  -4.20%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.50"
$ws.Range("E31").Value = "  +1.52%  "
$ws.Range("B32").Value = "dogwifhat"
$ws.Range("C32").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.07"
$ws.Range("E32").Value = "  +6.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.114"
$ws.Range("E33").Value = "  -1.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.17"
$ws.Range("E34").Value = "  +0.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.13"
$ws.Range("E35").Value = "  -5.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.612.51"
$ws.Range("E37").Value = "  +0.39%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.62"
$ws.Range("E38").Value = "  +4.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.393"
$ws.Range("E39").Value = "  -4.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.60"
$ws.Range("E40").Value = "  -3.99%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0770"
$ws.Range("E41").Value = "  -5.25%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "496.67"
$ws.Range("E42").Value = "  -8.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.136"
$ws.Range("E43").Value = "  -2.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0461"
$ws.Range("E44").Value = "  -3.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.90"
$ws.Range("E45").Value = "  -1.71%  "
$ws.Range("E46").Value = "  +2.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.33"
$ws.Range("E47").Value = "  -4.84%  "
$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  +0.42%  "
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.73"
$ws.Range("E49").Value = "  -6.95%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.01"
$ws.Range("E50").Value = "  -2.79%  "
$ws.Range("E51").Value = "  -8.45%  "
